# ProviderDataUploadTemplate.xlsx edit:
#   - Rename "Sheet1" -> "Vendor"
#   - Add a new "IssueList" sheet after "Vendor"
#   - Insert a "VendorCode" column at the front of the Vendor sheet
#   - Move the old IssueNo/IssueItem/Owner columns from Vendor to the new
#     IssueList sheet (tagged with the VendorCode) and add one new issue row

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Vendor"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "IssueList"

# --- Vendor sheet: insert a new "VendorCode" column at A -------------------
$ws1.Columns.Item(1).Insert()

# Give the new header cell the same look as the rest of row 1 (A1) and
# carry the data-row border style down into A2.
$ws1.Range("B1").Copy()
$ws1.Range("A1").PasteSpecial(-4122)
$ws1.Range("A1").Value = "VendorCode"
$ws1.Range("A2").Value = "ABC234"
$ws1.Columns.Item(1).ColumnWidth = 11.333333333333332

# --- IssueList sheet: relocate the old I:K (now J:L) issue columns ---------
$ws1.Range("J1:L2").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$ws1.Range("J1:L2").Copy()
$ws2.Range("B1").PasteSpecial(-4163)
$ws1.Range("J1:L2").ClearContents()

$ws2.Range("A1").Value = "VendorCode"
$ws2.Range("A2").Value = "ABC234"
$ws2.Range("A3").Value = "ABC234"

$ws2.Range("B3").Value = 3434
$ws2.Range("C3").Value = "dfdffd"
$ws2.Range("D3").Value = "Robert Sandford"

$ws2.Columns.Item(1).ColumnWidth = 11.333333333333332
$ws2.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws2.Columns.Item(3).ColumnWidth = 9.5
$ws2.Columns.Item(4).ColumnWidth = 14.666666666666666

# --- Selections matching the authored workbook ------------------------------
$ws1.Range("H2").Select()
$ws2.Range("D5").Select()
$ws1.Activate()
